# Auto-applied edit: weekly crime-data refresh (new week's figures).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-NumCell($addr, [double]$val, $donor) {
    $ws.Range($addr).Value = $val
    $ws.Range($donor).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

function Set-TextCell($addr, $text, $donor) {
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($donor).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# Row 15
Set-TextCell "D15" "0" "A14"
Set-TextCell "E15" "***.*" "A14"
Set-TextCell "F15" "0" "A14"
Set-NumCell "H15" -100 "H15"
Set-NumCell "N15" -62.5 "H15"

# Row 16
Set-NumCell "C16" 12 "G15"
Set-NumCell "D16" 14 "G15"
Set-NumCell "E16" -14.285714285714 "H15"
Set-NumCell "F16" 30 "G15"
Set-NumCell "G16" 63 "G15"
Set-NumCell "H16" -52.380952380952 "H15"
Set-NumCell "I16" 284 "G15"
Set-NumCell "J16" 374 "G15"
Set-NumCell "K16" -24.064171122994 "H15"
Set-NumCell "L16" 30.275229357798 "H15"
Set-NumCell "M16" 222.727272727273 "H15"
Set-NumCell "N16" -81.665590703679 "H15"

# Row 17
Set-NumCell "C17" 11 "G15"
Set-NumCell "D17" 10 "G15"
Set-NumCell "E17" 10 "H15"
Set-NumCell "F17" 54 "G15"
Set-NumCell "G17" 46 "G15"
Set-NumCell "H17" 17.391304347826 "H15"
Set-NumCell "I17" 299 "G15"
Set-NumCell "J17" 261 "G15"
Set-NumCell "K17" 14.559386973180 "H15"
Set-NumCell "L17" 4.912280701754 "H15"
Set-NumCell "M17" 157.758620689655 "H15"
Set-NumCell "N17" -21.315789473684 "H15"

# Row 18
Set-NumCell "D18" 12 "G15"
Set-NumCell "E18" -33.333333333333 "H15"
Set-NumCell "F18" 37 "G15"
Set-NumCell "G18" 53 "G15"
Set-NumCell "H18" -30.188679245283 "H15"
Set-NumCell "I18" 262 "G15"
Set-NumCell "J18" 412 "G15"
Set-NumCell "K18" -36.407766990291 "H15"
Set-NumCell "L18" 6.072874493927 "H15"
Set-NumCell "M18" 32.994923857868 "H15"
Set-NumCell "N18" -83.301465901848 "H15"

# Row 19
Set-NumCell "C19" 43 "G15"
Set-NumCell "D19" 51 "G15"
Set-NumCell "E19" -15.686274509803 "H15"
Set-NumCell "F19" 157 "G15"
Set-NumCell "G19" 183 "G15"
Set-NumCell "H19" -14.207650273224 "H15"
Set-NumCell "I19" 1413 "G15"
Set-NumCell "J19" 1316 "G15"
Set-NumCell "K19" 7.370820668693 "H15"
Set-NumCell "L19" 93.827160493827 "H15"
Set-NumCell "M19" 3.063457330415 "H15"
Set-NumCell "N19" -74.920127795527 "H15"

# Row 20
Set-TextCell "C20" "0" "A14"
Set-NumCell "D20" 2 "G15"
Set-NumCell "E20" -100 "H15"
Set-NumCell "F20" 7 "G15"
Set-NumCell "G20" 7 "G15"
Set-NumCell "H20" 0 "H15"
Set-NumCell "J20" 47 "G15"
Set-NumCell "K20" -19.148936170212 "H15"
Set-NumCell "L20" 31.034482758620 "H15"
Set-NumCell "N20" -82.159624413145 "H15"

# Row 21
Set-NumCell "C21" 74 "C21"
Set-NumCell "D21" 89 "C21"
Set-NumCell "E21" -16.853932584269 "E21"
Set-NumCell "F21" 285 "C21"
Set-NumCell "G21" 355 "C21"
Set-NumCell "H21" -19.718309859154 "E21"
Set-NumCell "I21" 2307 "C21"
Set-NumCell "J21" 2429 "C21"
Set-NumCell "K21" -5.022643062988 "E21"
Set-NumCell "L21" 51.676528599605 "E21"
Set-NumCell "M21" 28.810720268006 "E21"
Set-NumCell "N21" -75.394624573378 "E21"

# Row 22
Set-NumCell "C22" 7 "G15"
Set-NumCell "D22" 2 "G15"
Set-NumCell "E22" 250 "H15"
Set-NumCell "F22" 14 "G15"
Set-NumCell "G22" 12 "G15"
Set-NumCell "H22" 16.666666666666 "H15"
Set-NumCell "I22" 130 "G15"
Set-NumCell "J22" 111 "G15"
Set-NumCell "K22" 17.117117117117 "H15"
Set-NumCell "L22" 60.493827160493 "H15"
Set-NumCell "M22" 56.626506024096 "H15"

# Row 24
Set-NumCell "C24" 73 "G15"
Set-NumCell "D24" 62 "G15"
Set-NumCell "E24" 17.741935483871 "H15"
Set-NumCell "F24" 330 "G15"
Set-NumCell "G24" 295 "G15"
Set-NumCell "H24" 11.864406779661 "H15"
Set-NumCell "I24" 2411 "G15"
Set-NumCell "J24" 1896 "G15"
Set-NumCell "K24" 27.162447257384 "H15"
Set-NumCell "L24" 94.122383252818 "H15"
Set-NumCell "M24" -12.676566461427 "H15"

# Row 25
Set-NumCell "C25" 20 "G15"
Set-NumCell "D25" 11 "G15"
Set-NumCell "E25" 81.818181818181 "H15"
Set-NumCell "F25" 100 "G15"
Set-NumCell "G25" 66 "G15"
Set-NumCell "H25" 51.515151515151 "H15"
Set-NumCell "I25" 643 "G15"
Set-NumCell "J25" 530 "G15"
Set-NumCell "K25" 21.320754716981 "H15"
Set-NumCell "L25" 20.864661654135 "H15"
Set-NumCell "M25" 90.236686390532 "H15"

# Row 26
Set-TextCell "C26" "0" "A14"
Set-TextCell "D26" "0" "A14"
Set-TextCell "E26" "***.*" "A14"
Set-NumCell "F26" 1 "G15"
Set-NumCell "H26" -75 "H15"
Set-NumCell "L26" -29.166666666666 "H15"

# Row 27
Set-NumCell "C27" 3 "G15"
Set-NumCell "D27" 5 "G15"
Set-NumCell "E27" -40 "H15"
Set-NumCell "G27" 17 "G15"
Set-NumCell "H27" 11.764705882352 "H15"
Set-NumCell "I27" 135 "G15"
Set-NumCell "J27" 129 "G15"
Set-NumCell "K27" 4.651162790697 "H15"
Set-NumCell "L27" 64.634146341463 "H15"

# Row 30
Set-NumCell "D30" 1 "G15"
Set-NumCell "G30" 5 "G15"
Set-NumCell "J30" 17 "G15"
Set-NumCell "K30" -52.941176470588 "H15"

# Volume/Number header text: 'Number  30' -> 'Number  31'
$hdr = $ws.Range("A8")
$hdr.Characters(21, 2).Text = "31"

# Reporting week dates: 7/24/2023-7/30/2023 -> 7/31/2023-8/6/2023
$wk = $ws.Range("C9")
$wk.Characters(27, 9).Text = "7/31/2023"
$wk.Characters(47, 9).Text = "8/6/2023"
